# Worked on temporal resolution
# Update the "Demand" sheet (EU27.Elec time series): replace the single
# annual demand figure with a 12-step intertemporal series, and make the
# Demand tab the active one (it had been the SupIm tab before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand")

# Correct the existing value for step 1 and add the new steps 2-12, all
# carrying the same per-step demand figure.
$demandPerStep = 330708333
for ($step = 1; $step -le 12; $step++) {
    $row = $step + 2
    $ws.Cells.Item($row, 1).Value = $step
    $ws.Cells.Item($row, 2).Value = $demandPerStep
}

# Column B is now wide enough to need an explicit width so the numbers
# stay fully visible.
$ws.Columns.Item(2).ColumnWidth = 10.17

# Make "Demand" the active sheet/tab and leave the selection near the
# bottom of the data the author just entered.
$ws.Activate()
$ws.Range("C12").Select() | Out-Null
